$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 10000
$ws.Range("I51").Value = 10000
$ws.Range("K51").Value = 10000
$ws.Range("M51").Value = -9516
$ws.Range("H57").Value = 24432
$ws.Range("J57").Value = 24432
$ws.Range("L57").Value = 73296
$ws.Range("N57").Value = -74294
$ws.Range("H58").Value = 991.25
$ws.Range("I58").Value = 82.5
$ws.Range("J58").Value = 1900
$ws.Range("K58").Value = 247.5
$ws.Range("L58").Value = 5700
$ws.Range("M58").Value = -97.5
$ws.Range("N58").Value = -6000
$ws.Range("H88").Value = 837.75
$ws.Range("J88").Value = 1133
$ws.Range("L88").Value = 1133
$ws.Range("N88").Value = -1945
$ws.Range("H91").Value = 837.75
$ws.Range("J91").Value = 1133
$ws.Range("L91").Value = 1133
$ws.Range("N91").Value = -3941
$ws.Range("H99").Value = 209
$ws.Range("I99").Value = 191
$ws.Range("J99").Value = 245
$ws.Range("K99").Value = 573
$ws.Range("L99").Value = 735
$ws.Range("M99").Value = 925
$ws.Range("N99").Value = -3731
$ws.Range("H112").Value = 1049.8113
$ws.Range("J112").Value = 1059.6078
$ws.Range("L112").Value = 3178.8234
$ws.Range("N112").Value = -5394.8234
$ws.Range("H129").Value = 257121.89
$ws.Range("I129").Value = 425
$ws.Range("J129").Value = 303794.06
$ws.Range("K129").Value = 1275
$ws.Range("L129").Value = 911382.1799999999
$ws.Range("M129").Value = 3725
$ws.Range("N129").Value = -921382.1799999999
$ws.Range("H132").Value = 4121.5415
$ws.Range("I132").Value = 4500.524
$ws.Range("J132").Value = 1468.6666
$ws.Range("K132").Value = 13501.572
$ws.Range("L132").Value = 4405.9998
$ws.Range("M132").Value = -10971.572
$ws.Range("N132").Value = -9465.9998
$ws.Range("H134").Value = 47629.668
$ws.Range("J134").Value = 47629.668
$ws.Range("L134").Value = 47629.668
$ws.Range("N134").Value = -57769.668
$ws.Range("H138").Value = 2503.3948
$ws.Range("I138").Value = 1452.2222
$ws.Range("J138").Value = 3449.45
$ws.Range("K138").Value = 4356.6666
$ws.Range("L138").Value = 10348.35
$ws.Range("M138").Value = 783.3334000000004
$ws.Range("N138").Value = -20628.35

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5689.977
$ws.Range("I32").Value = 4403.6
$ws.Range("J32").Value = 10692.556
$ws.Range("K32").Value = 4403.6
$ws.Range("L32").Value = 10692.556
$ws.Range("M32").Value = -4116.6
$ws.Range("N32").Value = -11266.556
$ws.Range("H61").Value = 3049.8438
$ws.Range("I61").Value = 2807.2222
$ws.Range("J61").Value = 4360
$ws.Range("K61").Value = 2807.2222
$ws.Range("L61").Value = 4360
$ws.Range("M61").Value = -2595.2222
$ws.Range("N61").Value = -4784
$ws.Range("H132").Value = 18099.656
$ws.Range("I132").Value = 2339.65
$ws.Range("J132").Value = 44366.332
$ws.Range("K132").Value = 7018.950000000001
$ws.Range("L132").Value = 133098.996
$ws.Range("M132").Value = -4488.950000000001
$ws.Range("N132").Value = -138158.996
$ws.Range("H134").Value = 62000
$ws.Range("J134").Value = 62000
$ws.Range("L134").Value = 62000
$ws.Range("N134").Value = -72140
$ws.Range("H136").Value = 3049.8438
$ws.Range("I136").Value = 2807.2222
$ws.Range("J136").Value = 4360
$ws.Range("K136").Value = 8421.6666
$ws.Range("L136").Value = 13080
$ws.Range("M136").Value = -5871.6666
$ws.Range("N136").Value = -18180

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1776.7142
$ws.Range("I86").Value = 1608.409
$ws.Range("J86").Value = 2061.5386
$ws.Range("K86").Value = 1608.409
$ws.Range("L86").Value = 2061.5386
$ws.Range("M86").Value = -485.4090000000001
$ws.Range("N86").Value = -4307.5386
$ws.Range("H89").Value = 1776.7142
$ws.Range("I89").Value = 1608.409
$ws.Range("J89").Value = 2061.5386
$ws.Range("K89").Value = 8042.045
$ws.Range("L89").Value = 10307.693
$ws.Range("M89").Value = -2426.045
$ws.Range("N89").Value = -21539.693
$ws.Range("H94").Value = 673.3333
$ws.Range("I94").Value = 497.89474
$ws.Range("J94").Value = 1090
$ws.Range("K94").Value = 497.89474
$ws.Range("L94").Value = 1090
$ws.Range("M94").Value = -46.89474000000001
$ws.Range("N94").Value = -1992
$ws.Range("H134").Value = 3322.5557
$ws.Range("I134").Value = 3085.425
$ws.Range("J134").Value = 5219.6
$ws.Range("K134").Value = 9256.275000000001
$ws.Range("L134").Value = 15658.8
$ws.Range("M134").Value = -6721.275000000001
$ws.Range("N134").Value = -20728.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 190
$ws.Range("I22").Value = 190
$ws.Range("K22").Value = 190
$ws.Range("M22").Value = 160
$ws.Range("H31").Value = 3725.75
$ws.Range("I31").Value = 948.53845
$ws.Range("J31").Value = 5625.9473
$ws.Range("K31").Value = 948.53845
$ws.Range("L31").Value = 5625.9473
$ws.Range("M31").Value = -653.53845
$ws.Range("N31").Value = -6215.9473
$ws.Range("H34").Value = 3725.75
$ws.Range("I34").Value = 948.53845
$ws.Range("J34").Value = 5625.9473
$ws.Range("K34").Value = 948.53845
$ws.Range("L34").Value = 5625.9473
$ws.Range("M34").Value = -746.53845
$ws.Range("N34").Value = -6029.9473
$ws.Range("H58").Value = 20786
$ws.Range("I58").Value = 1556
$ws.Range("J58").Value = 64053.5
$ws.Range("K58").Value = 1556
$ws.Range("L58").Value = 64053.5
$ws.Range("M58").Value = -1353
$ws.Range("N58").Value = -64459.5
$ws.Range("H136").Value = 20786
$ws.Range("I136").Value = 1556
$ws.Range("J136").Value = 64053.5
$ws.Range("K136").Value = 4668
$ws.Range("L136").Value = 192160.5
$ws.Range("M136").Value = -2118
$ws.Range("N136").Value = -197260.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 13.181818
$ws.Range("I2").Value = 16.666666
$ws.Range("J2").Value = 9
$ws.Range("K2").Value = 99.999996
$ws.Range("L2").Value = 54
$ws.Range("M2").Value = 13.000004
$ws.Range("N2").Value = -280
$ws.Range("H5").Value = 1085.0264
$ws.Range("I5").Value = 930.0769
$ws.Range("J5").Value = 1420.75
$ws.Range("K5").Value = 2790.2307
$ws.Range("L5").Value = 4262.25
$ws.Range("M5").Value = -2678.2307
$ws.Range("N5").Value = -4486.25
$ws.Range("H31").Value = 980
$ws.Range("I31").Value = 1000
$ws.Range("J31").Value = 900
$ws.Range("K31").Value = 3000
$ws.Range("L31").Value = 2700
$ws.Range("M31").Value = -2712
$ws.Range("N31").Value = -3276
$ws.Range("H74").Value = 9874.75
$ws.Range("J74").Value = 9874.75
$ws.Range("L74").Value = 29624.25
$ws.Range("N74").Value = -31746.25
$ws.Range("H77").Value = 9874.75
$ws.Range("J77").Value = 9874.75
$ws.Range("L77").Value = 88872.75
$ws.Range("N77").Value = -99480.75
$ws.Range("H131").Value = 716.9091
$ws.Range("J131").Value = 723.30524
$ws.Range("L131").Value = 2169.91572
$ws.Range("N131").Value = -12249.91572
$ws.Range("H135").Value = 1085.0264
$ws.Range("I135").Value = 930.0769
$ws.Range("J135").Value = 1420.75
$ws.Range("K135").Value = 8370.6921
$ws.Range("L135").Value = 12786.75
$ws.Range("M135").Value = -5835.6921
$ws.Range("N135").Value = -17856.75
$ws.Range("H139").Value = 2090.5
$ws.Range("I139").Value = 1253.0741
$ws.Range("J139").Value = 3420.5293
$ws.Range("K139").Value = 3759.2223
$ws.Range("L139").Value = 10261.5879
$ws.Range("M139").Value = 1380.7777
$ws.Range("N139").Value = -20541.5879

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2383.0667
$ws.Range("I122").Value = 1014.7
$ws.Range("J122").Value = 5119.8
$ws.Range("K122").Value = 3044.1
$ws.Range("L122").Value = 15359.4
$ws.Range("M122").Value = -594.1000000000004
$ws.Range("N122").Value = -20259.4
$ws.Range("H126").Value = 2550.587
$ws.Range("I126").Value = 2105.739
$ws.Range("J126").Value = 2995.4348
$ws.Range("K126").Value = 6317.217000000001
$ws.Range("L126").Value = 8986.304400000001
$ws.Range("M126").Value = -3847.217000000001
$ws.Range("N126").Value = -13926.3044

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4010.353
$ws.Range("I7").Value = 4120
$ws.Range("J7").Value = 3853.7144
$ws.Range("K7").Value = 4120
$ws.Range("L7").Value = 3853.7144
$ws.Range("M7").Value = -4008
$ws.Range("N7").Value = -4077.7144
$ws.Range("H22").Value = 2588.5908
$ws.Range("I22").Value = 3034.4375
$ws.Range("K22").Value = 3034.4375
$ws.Range("M22").Value = -2739.4375
$ws.Range("H27").Value = 2588.5908
$ws.Range("I27").Value = 3034.4375
$ws.Range("K27").Value = 3034.4375
$ws.Range("M27").Value = -2927.4375
$ws.Range("H55").Value = 121.07692
$ws.Range("I55").Value = 85
$ws.Range("J55").Value = 178.8
$ws.Range("K55").Value = 85
$ws.Range("L55").Value = 178.8
$ws.Range("M55").Value = 88
$ws.Range("N55").Value = -524.8
$ws.Range("H122").Value = 787270.7
$ws.Range("I122").Value = 936503.2
$ws.Range("J122").Value = 3800
$ws.Range("K122").Value = 2809509.6
$ws.Range("L122").Value = 11400
$ws.Range("M122").Value = -2807059.6
$ws.Range("N122").Value = -16300
$ws.Range("H126").Value = 4010.353
$ws.Range("I126").Value = 4120
$ws.Range("J126").Value = 3853.7144
$ws.Range("K126").Value = 12360
$ws.Range("L126").Value = 11561.1432
$ws.Range("M126").Value = -9890
$ws.Range("N126").Value = -16501.1432
$ws.Range("H132").Value = 2167.2727
$ws.Range("I132").Value = 1326.8
$ws.Range("K132").Value = 3980.4
$ws.Range("M132").Value = -1450.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 3933
$ws.Range("J15").Value = 3933
$ws.Range("L15").Value = 3933
$ws.Range("N15").Value = -4509
$ws.Range("H21").Value = 3000
$ws.Range("J21").Value = 3000
$ws.Range("L21").Value = 3000
$ws.Range("N21").Value = -3470
$ws.Range("H35").Value = 3000
$ws.Range("J35").Value = 3000
$ws.Range("L35").Value = 3000
$ws.Range("N35").Value = -3580
$ws.Range("H126").Value = 1488.4375
$ws.Range("I126").Value = 1115.963
$ws.Range("K126").Value = 3347.889
$ws.Range("M126").Value = -877.8890000000001
$ws.Range("H132").Value = 1514.8214
$ws.Range("I132").Value = 1310.1333
$ws.Range("J132").Value = 1751
$ws.Range("K132").Value = 3930.3999
$ws.Range("L132").Value = 5253
$ws.Range("M132").Value = -1400.3999
$ws.Range("N132").Value = -10313
$ws.Range("H140").Value = 45199.75
$ws.Range("J140").Value = 45199.75
$ws.Range("L140").Value = 45199.75
$ws.Range("N140").Value = -55559.75
